$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'ZARA India | New Collection Online'
$ws.Cells.Item(2, 2).Value = 'https://www.zara.com/in/'
$ws.Cells.Item(2, 3).Value = 'WebWe are present in over 156 stores, please select yours. here. Latest trends in clothing for women, men & kids at ZARA online. Find new arrivals, fashion catalogs, collections & …'
$ws.Cells.Item(2, 4).Value = 'N/A'

$ws.Cells.Item(3, 1).Value = 'ZARA Official Website'
$ws.Cells.Item(3, 2).Value = 'https://www.zara.com/'
$ws.Cells.Item(3, 3).Value = 'WebZARA Official Website. Presione enter para accesibilidad para personas ciegas. Presione enter para navegar por el teclado. Presione enter para acceder al menú de accesibilidad. …'
$ws.Cells.Item(3, 4).Value = 'N/A'

$ws.Cells.Item(4, 1).Value = 'Women''s New In Clothes | Explore our New Arrivals | ZARA India'
$ws.Cells.Item(4, 2).Value = 'https://www.zara.com/in/en/woman-new-in-l1180.html'
$ws.Cells.Item(4, 3).Value = 'Webzara origins NEW; edition NEW; #selectedby NEW /// new; jackets; puffers; trousers; jeans; sweaters | cardigans; hoodies | sweatshirts; t-shirts; overshirts; shirts; polo shirts; linen; …'
$ws.Cells.Item(4, 4).Value = 'N/A'

$ws.Cells.Item(5, 1).Value = 'Women''s Jackets | ZARA India'
$ws.Cells.Item(5, 2).Value = 'https://www.zara.com/in/en/woman-jackets-l1114.html'
$ws.Cells.Item(5, 3).Value = 'WebZARA''s women''s jackets and waistcoats are easy to wear and complement any look. Our women''s corduroy jackets will keep the cold at bay and add instant edge to any outfit, …'
$ws.Cells.Item(5, 4).Value = 'N/A'

$ws.Cells.Item(6, 1).Value = 'LEATHER JACKET - Dark brown | ZARA India'
$ws.Cells.Item(6, 2).Value = 'https://www.zara.com/in/en/leather-jacket-p02521105.html'
$ws.Cells.Item(6, 3).Value = 'WebCropped fit jacket made of leather treated with a process that gives it a creased appearance. Lapel collar and long sleeves. Hip welt pockets and interior pocket detail. Front zip …'
$ws.Cells.Item(6, 4).Value = 'N/A'

$ws.Cells.Item(7, 1).Value = 'Women''s Handbags | Explore our New Arrivals | ZARA India'
$ws.Cells.Item(7, 2).Value = 'https://www.zara.com/in/en/woman-bags-handbags-l1037.html'
$ws.Cells.Item(7, 3).Value = 'WebWomen’s handbags come in all shapes for this season, from the traditional to the unexpected. Vintage silhouettes such as Nineties shoulder bags and Fifties top handle …'
$ws.Cells.Item(7, 4).Value = 'N/A'

$ws.Cells.Item(8, 1).Value = 'Dresses for Women | ZARA India'
$ws.Cells.Item(8, 2).Value = 'https://www.zara.com/in/en/woman-dresses-l1066.html'
$ws.Cells.Item(8, 3).Value = 'WebZARA''s short dresses include both classic and on-trend styles. The short black dress, also known as little black dress, is an emblematic garment basic for any woman''s wardrobe. One of its many qualities is not only the great potential in terms of versatility; it is a piece that exudes unpretentious elegance wherever it goes.'
$ws.Cells.Item(8, 4).Value = 'N/A'

$ws.Cells.Item(9, 1).Value = 'Women''s Clothes | ZARA United Kingdom'
$ws.Cells.Item(9, 2).Value = 'https://www.zara.com/uk/en/woman-mkt1000.html'
$ws.Cells.Item(9, 3).Value = 'WebWeekly new trends in clothes, shoes & accessories at ZARA online. FREE SHIPPING* for you to try on at your leisure.'
$ws.Cells.Item(9, 4).Value = 'N/A'

$ws.Cells.Item(10, 1).Value = 'Men´s Clothes | ZARA India'
$ws.Cells.Item(10, 2).Value = 'https://www.zara.com/in/en/man-mkt534.html'
$ws.Cells.Item(10, 3).Value = 'WebThe latest clothes, shoes and accessories for men every week at ZARA online. Enter now and discover all the shirts of the new collection'
$ws.Cells.Item(10, 4).Value = 'N/A'

$ws.Cells.Item(11, 1).Value = 'Women''s Just In Clothes | Explore our New Arrivals - ZARA'
$ws.Cells.Item(11, 2).Value = 'https://www.zara.com/us/en/woman-new-in-l1180.html'
$ws.Cells.Item(11, 3).Value = 'WebFREE SHIPPING. New clothes and accessories updated weekly at ZARA online. Stay in style with seasonal trends.'
$ws.Cells.Item(11, 4).Value = 'N/A'

$ws.Cells.Item(12, 1).Value = 'ZARA Canada | New Collection Online'
$ws.Cells.Item(12, 2).Value = 'https://www.zara.com/ca/'
$ws.Cells.Item(12, 3).Value = 'WebDiscover the new ZARA collection online. The latest trends for Woman, Man, Kids and next season’s ad campaigns.'
$ws.Cells.Item(12, 4).Value = 'N/A'

$ws.Cells.Item(13, 1).Value = 'ZARA Philippines | New Collection Online'
$ws.Cells.Item(13, 2).Value = 'https://www.zara.com/ph/'
$ws.Cells.Item(13, 3).Value = 'WebLatest trends in clothing for women, men & kids at ZARA online. Find new arrivals, fashion catalogs, collections & lookbooks every week.'
$ws.Cells.Item(13, 4).Value = 'N/A'

$ws.Cells.Item(14, 1).Value = 'Dresses for Women | ZARA United States'
$ws.Cells.Item(14, 2).Value = 'https://www.zara.com/us/en/woman-dresses-l1066.html'
$ws.Cells.Item(14, 3).Value = 'WebZara''s women''s dresses collection comprises a wide variety of styles. From the latest arrivals to pieces for a basic wardrobe, both on-trend and timeless designs are featured in the collection. A black dress represents the highest level of elegance. Zara''s selection includes different alternatives to the iconic little black dress, a key piece ...'
$ws.Cells.Item(14, 4).Value = 'N/A'

$ws.Cells.Item(15, 1).Value = 'Men´s New In Clothes | Explore our New Arrivals | ZARA India'
$ws.Cells.Item(15, 2).Value = 'https://www.zara.com/in/en/man-new-in-l711.html'
$ws.Cells.Item(15, 3).Value = 'WebOur edit of new men''s clothes takes in timeless wardrobe heroes and new trends in fashion. From staples - including T-shirts, vests, knitwear, sleepwear and underwear - to smart occasionwear, such as shirts, blazers and co-ord suits, the complete closet has been considered and is waiting to shop online.'
$ws.Cells.Item(15, 4).Value = 'N/A'

$ws.Cells.Item(16, 1).Value = 'Women''s Special Prices | Explore our New Arrivals | ZARA India'
$ws.Cells.Item(16, 2).Value = 'https://www.zara.com/in/en/woman-special-prices-l1314.html'
$ws.Cells.Item(16, 3).Value = 'WebThe Zara special prices edit is perfect for finding that one thing your wardrobe is missing, whether it be a little black dress, blazer, a pair of jeans, or staples like vests, T-shirts and lingerie. Available in store and online, shop clothes, accessories and footwear on sale. TURTLENECK SWEATER +3 ₹ 2,290.00-30 % ₹ 1,590.00. TURTLENECK SWEATER …'
$ws.Cells.Item(16, 4).Value = 'N/A'

$ws.Cells.Item(17, 1).Value = 'ZARA UAE - Dubai/Sharjah/Ajman/UAQ/Fujairah | New …'
$ws.Cells.Item(17, 2).Value = 'https://www.zara.com/ae/'
$ws.Cells.Item(17, 3).Value = 'WebLatest trends in clothing for women, men & kids at ZARA online. Find new arrivals, fashion catalogs, collections & lookbooks every week.'
$ws.Cells.Item(17, 4).Value = 'N/A'

$ws.Cells.Item(18, 1).Value = 'Women''s Trousers | ZARA India'
$ws.Cells.Item(18, 2).Value = 'https://www.zara.com/in/en/woman-trousers-l1335.html'
$ws.Cells.Item(18, 3).Value = 'WebZARA''s collection of women''s trousers connects every woman to her ideal fit. The wide variety of cuts and styles allows you to create a versatile wardrobe with trousers for every occasion. Women''s dress trousers are one of the most flattering garments: they are slimming and offer the possibility to select the right fit for each body type.'
$ws.Cells.Item(18, 4).Value = 'N/A'

$ws.Cells.Item(19, 1).Value = 'ZARA Ireland | New Collection Online'
$ws.Cells.Item(19, 2).Value = 'https://www.zara.com/ie/'
$ws.Cells.Item(19, 3).Value = 'WebLatest trends in clothing for women, men & kids at ZARA online. Find new arrivals, fashion catalogs, collections & lookbooks every week.'
$ws.Cells.Item(19, 4).Value = 'N/A'

$ws.Cells.Item(20, 1).Value = 'ZARA Hungary / Hungary | New Collection Online'
$ws.Cells.Item(20, 2).Value = 'https://www.zara.com/hu/'
$ws.Cells.Item(20, 3).Value = 'WebLatest trends in clothing for women, men & kids at ZARA online. Find new arrivals, fashion catalogs, collections & lookbooks every week.'
$ws.Cells.Item(20, 4).Value = 'N/A'

$ws.Cells.Item(21, 1).Value = 'Women''s Jackets | ZARA India'
$ws.Cells.Item(21, 2).Value = 'https://www.zara.com/in/en/woman-jackets-l1114.html'
$ws.Cells.Item(21, 3).Value = 'WebZARA''s women''s jackets and waistcoats are easy to wear and complement any look. Our women''s corduroy jackets will keep the cold at bay and add instant edge to any outfit, thanks to their soft ribbed fabric which lasts in every wardrobe through the seasons. Meanwhile, our iconic women''s leather jackets, regardless of their fit, are timeless ...'
$ws.Cells.Item(21, 4).Value = 'N/A'

$ws.Cells.Item(22, 1).Value = 'ZARA Saudi Arabia | New Collection Online'
$ws.Cells.Item(22, 2).Value = 'https://www.zara.com/sa/en/'
$ws.Cells.Item(22, 3).Value = 'WebYes, continue on Saudi Arabia. No, go to the website for United States. We are present in over 156 stores, please select yours. here. Latest trends in clothing for women, men & kids at ZARA online. Find new arrivals, fashion catalogs, collections & lookbooks every week.'
$ws.Cells.Item(22, 4).Value = 'N/A'

$ws.Cells.Item(23, 1).Value = 'Women''s New In Clothes | Explore our New Arrivals - ZARA'
$ws.Cells.Item(23, 2).Value = 'https://www.zara.com/uk/en/woman-new-in-l1180.html'
$ws.Cells.Item(23, 3).Value = 'WebNew clothes and accessories updated weekly at ZARA online. Stay in style with seasonal trends.'
$ws.Cells.Item(23, 4).Value = 'N/A'

$ws.Cells.Item(24, 1).Value = 'Women''s Blazers | ZARA India'
$ws.Cells.Item(24, 2).Value = 'https://www.zara.com/in/en/woman-blazers-l1055.html'
$ws.Cells.Item(24, 3).Value = 'WebZara offers a collection of women''s blazers that focuses on the original silhouette of this jacket, including updated proposals as well as classic designs that never go out of fashion. For a sophisticated, seamless look, the black blazer and the white blazer have a special place. As contemporary alternatives, the cropped and oversize blazer carry the …'
$ws.Cells.Item(24, 4).Value = 'N/A'

$ws.Cells.Item(25, 1).Value = 'Women''s T-shirts | ZARA India'
$ws.Cells.Item(25, 2).Value = 'https://www.zara.com/in/en/woman-tshirts-l1362.html'
$ws.Cells.Item(25, 3).Value = 'WebThis collection of women''s T-shirts from Zara pays tribute to an everyday garment that holds a special place in women''s wardrobes. Black and white T-shirts are neutral and adaptable, which makes them suitable for a wide range of styles, while oversized and cropped T-shirts add an urban touch to casual looks. In addition, this garment is also suitable for delicate …'
$ws.Cells.Item(25, 4).Value = 'N/A'

$ws.Cells.Item(26, 1).Value = 'ZARA North Macedonia | New Collection Online'
$ws.Cells.Item(26, 2).Value = 'https://www.zara.com/mk/en/'
$ws.Cells.Item(26, 3).Value = 'WebHello, Yes, continue on North Macedonia. Yes, continue on North Macedonia. No, go to the website for United States. We are present in over 156 stores, please select yours. here. Latest trends in clothing for women, men & kids at ZARA online. Find new arrivals, fashion catalogs, collections & lookbooks every week.'
$ws.Cells.Item(26, 4).Value = 'N/A'

$ws.Cells.Item(27, 1).Value = 'Women''s Tops | ZARA India'
$ws.Cells.Item(27, 2).Value = 'https://www.zara.com/in/en/woman-tops-l1322.html'
$ws.Cells.Item(27, 3).Value = 'WebTops For Women. Our collection of women’s tops has something to balance with every outfit, whether you prefer off-the-shoulder styles, halternecks, camisoles or knitted jersey t-shirts. Classic black, white, pink and navy block colours sit alongside everything from sequin tops to statement animal prints, to original slogan T-shirts this season.'
$ws.Cells.Item(27, 4).Value = 'N/A'

$ws.Range("A28:D38").EntireRow.Delete()
